$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while assigning, so numeric-looking
# strings (e.g. "1.00", "19.20") are preserved verbatim as text,
# matching the inline-string cell type in the source workbook.
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "42.462.14"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "2.286.37"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "156.89"
$ws.Range("E5").Value = "  +15,572.23%  "
$ws.Range("D6").Value = "307.09"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").Value = "96.90"
$ws.Range("E7").Value = "  +6.19%  "
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").Value = "36.36"
$ws.Range("E11").Value = "  +13.59%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "2.640.85"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "14.63"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "2.294.26"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "0.805"
$ws.Range("E18").Value = "  +6.16%  "
$ws.Range("D19").Value = "42.361.78"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "12.87"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").Value = "67.86"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").Value = "242.81"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").Value = "36.65"
$ws.Range("E29").Value = "  +7.62%  "
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "161.50"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").Value = "17.39"
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("E39").Value = "  +5.32%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +6.42%  "
$ws.Range("D43").Value = "2.38"
$ws.Range("E43").Value = "  +16.88%  "
$ws.Range("D44").Value = "2.002.28"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "19.20"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "10.28"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "3.03"
$ws.Range("E48").Value = "  +6.31%  "
$ws.Range("D49").Value = "53.92"
$ws.Range("E49").Value = "  +4.76%  "
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").Value = "72.58"
$ws.Range("E51").Value = "  +0.01%  "

# Remove the temporary text-format override so styling matches the original
# (cells keep their plain/default style, no explicit numFmt left behind).
$dCol.ClearFormats()
